$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 5 with the new test case CP_GESCLSERDOM_004
$ws.Range("A5").Value = "CP_GESCLSERDOM_004"
$ws.Range("B5").Value = "Reconfiguración `ndel cliente"
$ws.Range("C5").Value = "Positivo"
$ws.Range("D5").Value = "eCenter"
$ws.Range("E5").Value = "Cliente seleccionado y visible en la tabla"
$ws.Range("F5").Value = "1. Clic en Opciones.`n2. Seleccionar Reconfiguración.`n3. Clic en Reconfigurar.`n4. Confirmar en el modal con Sí."
$ws.Range("G5").Value = "ID de cliente válido"
$ws.Range("H5").Value = "Se inicia el proceso de reconfiguración y se muestran barras de progreso."
$ws.Range("I5").Value = "Se visualiza correctamente el proceso de reconfiguracion"
$ws.Range("J5").Value = "OK"
$ws.Range("K5").Value = "SI"
$ws.Range("L5").Value = "N/A"

# Move the active selection from B4 to B5
$ws.Range("B5").Select()
